$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.883886666666666
$ws.Range("H2").Value = 5.65166
$ws.Range("I2").Value = 0.7298568945019562
$ws.Range("J2").Value = 0.7298568945019563
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.425703666666666
$ws.Range("N2").Value = 4.277111
$ws.Range("O2").Value = 0.04715501820393346
$ws.Range("P2").Value = 0.04715501820393346
$ws.Range("Q2").Value = 2.685864128251111
$ws.Range("R2").Value = 24.17277715426
$ws.Range("S2").Value = 0.03441641514650609
$ws.Range("T2").Value = 0.0344164151465061

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.883886666666666
$ws.Range("H3").Value = 5.65166
$ws.Range("I3").Value = 0.7298568945019562
$ws.Range("J3").Value = 0.7298568945019563
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.61433933333333
$ws.Range("N3").Value = 61.843018
$ws.Range("O3").Value = 0.6818173855147049
$ws.Range("P3").Value = 0.6818173855147048
$ws.Range("Q3").Value = 38.83507901220889
$ws.Range("R3").Value = 349.51571110988
$ws.Range("S3").Value = 0.4976291196092056
$ws.Range("T3").Value = 0.4976291196092056

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.883886666666666
$ws.Range("H4").Value = 5.65166
$ws.Range("I4").Value = 0.7298568945019562
$ws.Range("J4").Value = 0.7298568945019563
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.194356666666666
$ws.Range("N4").Value = 24.58307
$ws.Range("O4").Value = 0.2710275962813615
$ws.Range("P4").Value = 0.2710275962813615
$ws.Range("Q4").Value = 15.43723926624444
$ws.Range("R4").Value = 138.9351533962
$ws.Range("S4").Value = 0.1978113597462445
$ws.Range("T4").Value = 0.1978113597462445

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.353459
$ws.Range("H5").Value = 1.060377
$ws.Range("I5").Value = 0.1369373713601492
$ws.Range("J5").Value = 0.1369373713601492
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.425703666666666
$ws.Range("N5").Value = 4.277111
$ws.Range("O5").Value = 0.04715501820393346
$ws.Range("P5").Value = 0.04715501820393346
$ws.Range("Q5").Value = 0.5039277923163333
$ws.Range("R5").Value = 4.535350130847
$ws.Range("S5").Value = 0.006457284239286633
$ws.Range("T5").Value = 0.006457284239286635

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.353459
$ws.Range("H6").Value = 1.060377
$ws.Range("I6").Value = 0.1369373713601492
$ws.Range("J6").Value = 0.1369373713601492
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 20.61433933333333
$ws.Range("N6").Value = 61.843018
$ws.Range("O6").Value = 0.6818173855147049
$ws.Range("P6").Value = 0.6818173855147048
$ws.Range("Q6").Value = 7.286323766420667
$ws.Range("R6").Value = 65.57691389778601
$ws.Range("S6").Value = 0.09336628052003318
$ws.Range("T6").Value = 0.09336628052003318

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.353459
$ws.Range("H7").Value = 1.060377
$ws.Range("I7").Value = 0.1369373713601492
$ws.Range("J7").Value = 0.1369373713601492
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.194356666666666
$ws.Range("N7").Value = 24.58307
$ws.Range("O7").Value = 0.2710275962813615
$ws.Range("P7").Value = 0.2710275962813615
$ws.Range("Q7").Value = 2.896369113043333
$ws.Range("R7").Value = 26.06732201739
$ws.Range("S7").Value = 0.0371138066008294
$ws.Range("T7").Value = 0.03711380660082941

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.343827
$ws.Range("H8").Value = 1.031481
$ws.Range("I8").Value = 0.1332057341378944
$ws.Range("J8").Value = 0.1332057341378944
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.425703666666666
$ws.Range("N8").Value = 4.277111
$ws.Range("O8").Value = 0.04715501820393346
$ws.Range("P8").Value = 0.04715501820393346
$ws.Range("Q8").Value = 0.490195414599
$ws.Range("R8").Value = 4.411758731391
$ws.Range("S8").Value = 0.006281318818140733
$ws.Range("T8").Value = 0.006281318818140733

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.343827
$ws.Range("H9").Value = 1.031481
$ws.Range("I9").Value = 0.1332057341378944
$ws.Range("J9").Value = 0.1332057341378944
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.61433933333333
$ws.Range("N9").Value = 61.843018
$ws.Range("O9").Value = 0.6818173855147049
$ws.Range("P9").Value = 0.6818173855147048
$ws.Range("Q9").Value = 7.087766449962001
$ws.Range("R9").Value = 63.78989804965801
$ws.Range("S9").Value = 0.09082198538546607
$ws.Range("T9").Value = 0.09082198538546606

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.343827
$ws.Range("H10").Value = 1.031481
$ws.Range("I10").Value = 0.1332057341378944
$ws.Range("J10").Value = 0.1332057341378944
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.194356666666666
$ws.Range("N10").Value = 24.58307
$ws.Range("O10").Value = 0.2710275962813615
$ws.Range("P10").Value = 0.2710275962813615
$ws.Range("Q10").Value = 2.81744106963
$ws.Range("R10").Value = 25.35696962667
$ws.Range("S10").Value = 0.03610242993428763
$ws.Range("T10").Value = 0.03610242993428763

